$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 154-155, pushing existing rows 154..176 down to 156..178
$ws.Range("A154:A155").EntireRow.Insert()

# --- New row 154 ---
$ws.Range("A154").Value = 9
$ws.Range("B154").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C154").Value = "Metropolitana"
$ws.Range("D154").Value = 44522
$ws.Range("E154").Value = 13
$ws.Range("F154").Value = 100112043
$ws.Range("G154").Value = "Pepino ensalada"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 160
$ws.Range("K154").Value = 5000
$ws.Range("L154").Value = 6000
$ws.Range("M154").Value = 5500
$ws.Range("N154").Value = "$/caja 50 unidades"
$ws.Range("O154").Value = "Región de Arica y Parinacota"
$ws.Range("P154").Value = 110
$ws.Range("Q154").Value = 50
$ws.Range("R154").Value = "Hortaliza"

# --- New row 155 ---
$ws.Range("A155").Value = 9
$ws.Range("B155").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C155").Value = "Metropolitana"
$ws.Range("D155").Value = 44522
$ws.Range("E155").Value = 13
$ws.Range("F155").Value = 100112043
$ws.Range("G155").Value = "Pepino ensalada"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Segunda"
$ws.Range("J155").Value = 97
$ws.Range("K155").Value = 4000
$ws.Range("L155").Value = 4000
$ws.Range("M155").Value = 4000
$ws.Range("N155").Value = "$/caja 100 unidades"
$ws.Range("O155").Value = "Región de Arica y Parinacota"
$ws.Range("P155").Value = 40
$ws.Range("Q155").Value = 100
$ws.Range("R155").Value = "Hortaliza"
